# Added Test Data for UK Market
#
# The "Poland" worksheet is the closest existing template (same column
# widths, same row-5 height, same merged header cells) for a new market
# sheet, so duplicate it, rename the copy to "UK", and then adjust its
# content:
#   - insert the two "P32AR"/"P32DR" rows that Poland's list is missing
#     (UK's repeater list matches the full list, e.g. like Germany's)
#   - set the market name (B2) and ticket/reference code (B4)
#   - leave the selection on B2, as in the authored sheet

$wb = $excel.ActiveWorkbook

$poland = $wb.Worksheets.Item("Poland")
$poland.Copy($null, $poland) | Out-Null

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "UK"

# Poland's repeater list skips "P32AR" / "P32DR" - insert two rows above
# the current row 16 ("PR1DS") to make room for them.
$newSheet.Rows.Item(16).Resize(2).Insert() | Out-Null

# The blank inserted rows don't inherit the list's cell style, so copy it
# down from the row just above (A15, "MZXDR240") before filling them in.
$newSheet.Range("A15").Copy() | Out-Null
$newSheet.Range("A16:A17").PasteSpecial(-4122) | Out-Null

$newSheet.Range("A16").Value = "P32AR"
$newSheet.Range("A17").Value = "P32DR"

# Fill in B4 before B2 so the shared-string table gets the two new
# strings in the same order as the authored workbook.
$newSheet.Range("B4").Value = "NGC-2741/T3355"
$newSheet.Range("B2").Value = "UK Market"

$newSheet.Range("B2").Select() | Out-Null
